$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Top-of-table single-value rows (1-indexed) ---
# Row 1: 99.94 -> 0M
$t.Rows.Item(1).Cells.Item(1).Range.Text = "0M"
# Row 2: 0.4 -> 0M
$t.Rows.Item(2).Cells.Item(1).Range.Text = "0M"
# Row 3: 681 -> 0M
$t.Rows.Item(3).Cells.Item(1).Range.Text = "0M"
# Row 4: 1073 -> 2229
$t.Rows.Item(4).Cells.Item(1).Range.Text = "2229"
# Row 5: 0.00004 -> 0.00002
$t.Rows.Item(5).Cells.Item(1).Range.Text = "0.00002"

# Rows 6,7,8 (0.00052 / 0.00012 / 0.00003) collapse into a single row
# holding 0.00629; remove the now-redundant rows 7 and 8.
$t.Rows.Item(6).Cells.Item(1).Range.Text = "0.00629"
$t.Rows.Item(7).Delete()
$t.Rows.Item(7).Delete()

# Row 9 (was 0.00015) is untouched.
# Row 8 (was 0.00016, after the collapse above) -> 0.00007
$t.Rows.Item(8).Cells.Item(1).Range.Text = "0.00007"
# Row 9 (was 0.00018) -> 0.00028
$t.Rows.Item(9).Cells.Item(1).Range.Text = "0.00028"

# Row 10 (was 0.13231) expands into three rows: 0.00032 / 0.00035 / 0.39927
$t.Rows.Item(10).Cells.Item(1).Range.Text = "0.00032"
$newRow1 = $t.Rows.Add($t.Rows.Item(11))
$newRow1.Cells.Item(1).Range.Text = "0.00035"
$newRow2 = $t.Rows.Add($t.Rows.Item(12))
$newRow2.Cells.Item(1).Range.Text = "0.39927"

# --- Bottom-of-table rows that previously held tab-separated run sequences ---
# collapse down to the single values that moved up to rows 1-3 above.
$t.Rows.Item(44).Cells.Item(1).Range.Text = "99.94"
$t.Rows.Item(45).Cells.Item(1).Range.Text = "0.4"
$t.Rows.Item(46).Cells.Item(1).Range.Text = "681"
